# Update the "想去人数" (want-to-go count) figures in column F for the
# rows that changed between crawls, on both the "展览" and "全部类型"
# worksheets (sheet2 "演出" and sheet3 "本地生活" have no data rows and
# are unaffected).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 204
    3  = 1072
    6  = 4599
    8  = 383
    9  = 1352
    10 = 893
    12 = 986
    14 = 556
    16 = 257
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
